# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# The "Rules" sheet has a small rule table in B7:E11 (rule name / from / to /
# greeting). Row 11 ("R40") is being re-labelled "1" - a digit-looking label,
# so it must stay a *text* entry (not be silently re-interpreted as the
# number 1). Prefixing with a leading apostrophe is the standard Excel way
# of forcing a numeric-looking entry to be stored as text, matching the
# shared-string cell that shows up in the workbook's XML.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
